# Applies the edit described by the diff: a new data row is inserted at
# worksheet row 13 (pushing the former rows 13..121 down to 14..122), and
# the new row 13 is populated with its own record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; this shifts rows 13:121 down to 14:122
# and carries formatting (e.g. the date style on column D) along with it.
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new record's values.
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "Macroferia Regional de Talca"
$ws.Range("C13").Value = "Maule"
$ws.Range("D13").Value = 44921
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100103
$ws.Range("H13").Value = "Frutos de hueso (carozo)"
$ws.Range("I13").Value = 100103002
$ws.Range("J13").Value = "Ciruela"
$ws.Range("K13").Value = "Black Amber"
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 150
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("Q13").Value = "$/bandeja 18 kilos granel"
$ws.Range("R13").Value = "Región de O'Higgins"
$ws.Range("S13").Value = 667
$ws.Range("T13").Value = 18
